# Q3 Update - 2025
#
# 1) The placeholder "coo_name" text (column B) used by every data row changes
#    from "Bau54O" to "7y1bbS".
# 2) Row 537 gets corrected refugees/asylum_seekers/returned_refugees figures.
# 3) Row 545 gets a corrected asylum_seekers figure.
#
# All of the touched numeric-looking figures (refugees/asylum_seekers/
# returned_refugees) are stored as *text* in this workbook (General number
# format), so a direct `.Value = "499"` assignment would make Excel silently
# reinterpret the text as a number. To keep the cells text-typed - matching
# the original authoring - the new text is first written into a scratch cell
# that has been forced to Text format, then only the *value* (not the number
# format/style) is pasted into the destination cell, leaving the destination
# cell's original style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fromCSV")

# 1) Column B ("coo_name") placeholder text used by every data row.
#    These aren't numeric-looking, so a plain assignment keeps them as text.
$ws.Range("B2:B553").Value = "7y1bbS"

# 2) & 3) Numeric-looking text fields - write via a Text-formatted scratch
#    cell so the destination cells keep their original style/number format
#    and stay string-typed.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$scratch.Value = "22407"
$scratch.Copy()
$ws.Range("N537").PasteSpecial(-4163)   # xlPasteValues

$scratch.Value = "499"
$scratch.Copy()
$ws.Range("O537").PasteSpecial(-4163)

$scratch.Value = "126"
$scratch.Copy()
$ws.Range("P537").PasteSpecial(-4163)

$scratch.Value = "140"
$scratch.Copy()
$ws.Range("O545").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
